$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the shared-string-backed header/data labels ---------------
# "Account No" -> "Account Number", "Mobile No" -> "Phine Number",
# "Account Balance" -> "Balance", "Last Transaction Date/Time" -> "Last Trans Details"
# are realised by writing the new header text into row 2 below.

# --- Fill in the (previously blank) header row 2 with real labels -----
$ws.Range("A2").Value = "Name"
$ws.Range("B2").Value = "Account Number"
$ws.Range("C2").Value = "Phine Number"
$ws.Range("D2").Value = "Email Address"
$ws.Range("E2").Value = "Password"
$ws.Range("F2").Value = "Account Type"
$ws.Range("G2").Value = "Balance"
$ws.Range("H2").Value = "Last Deposite"
$ws.Range("I2").Value = "Last Transaction"
$ws.Range("J2").Value = "Last Trans Details"

# --- Add the new "Last Transaction" detail values for the data rows ---
$ws.Range("I3").Value = "+ 1500"
$ws.Range("I4").Value = "+ 2500"

# --- Make data rows 3-4 match the formatting of the header rows -------
$ws.Range("A1:J1").Copy()
$ws.Range("A3:J4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(3).RowHeight = 18
$ws.Rows.Item(4).RowHeight = 18
